# IWP TestData, Demo Verification Script: Committing after Windows 11 update
#
# Updates the latest-run Result/Date columns on the three bootstrap sheets
# to record a new (failing) execution pass.

$wb = $excel.ActiveWorkbook

# --- CreateModifyDeleteProfile ---------------------------------------
$ws1 = $wb.Worksheets.Item("CreateModifyDeleteProfile")
$ws1.Range("A2").Value = "Fail"
$ws1.Range("B2").Value = "Mon Aug 04 20:51:46 IST 2025"

# --- AddModifyDeleteCC -------------------------------------------------
$ws2 = $wb.Worksheets.Item("AddModifyDeleteCC")
$ws2.Range("A2").Value = "Fail"
$ws2.Range("B2").Value = "Mon Aug 04 20:49:27 IST 2025"

# --- AddModifyDeleteACH -------------------------------------------------
$ws3 = $wb.Worksheets.Item("AddModifyDeleteACH")
$ws3.Range("A2").Value = "Fail"
$ws3.Range("B2").Value = "Mon Aug 04 20:42:39 IST 2025"
$ws3.Range("A3").Value = "Fail"
$ws3.Range("B3").Value = "Mon Aug 04 20:44:03 IST 2025"
$ws3.Range("A4").Value = "Fail"
$ws3.Range("B4").Value = "Mon Aug 04 20:45:13 IST 2025"
